$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder header labels (row 1) to group columns by room type
# (kitchens, bedrooms, living_rooms) instead of the previous order.
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "kitchens_2"
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("D1").Value = "bedrooms_2"
$ws.Range("E1").Value = "living_rooms_1"
$ws.Range("F1").Value = "living_rooms_2"

# Rewrite the one-hot data rows (rows 2-7) to match the new column order.
$data = @(
    @(0,0,0,1,0,0),
    @(0,1,0,0,0,0),
    @(0,0,0,0,0,1),
    @(0,0,1,0,0,0),
    @(1,0,0,0,0,0),
    @(0,0,0,0,1,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
